# Arquitectura de Redes - agregar manejo automatico de numero de actividades
# Se actualizan Docente/Materia/Carrera/Semestre (normalizados a mayusculas
# donde corresponde) y se completan las filas de Tema / Trabajo Independiente
# para las 18 actividades del semestre.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docente  = "Ing. Carlos Guzman"
$materia  = "ARQUITECTURA DE REDES"
$carrera  = "REDES Y TELECOMUNICACIONES"
$semestre = "SEGUNDO"

# Actividad (E), Tema (F), Trabajo Independiente (G) por fila.
$actividades = @(
    @(1,  "Encuadre",            "Firmar Encuadre"),
    @(2,  "Prueba Diagnóstico",  "SN"),
    @(3,  "Introducción",        "SN"),
    @(4,  "Contenidos Varios",   "Tipos de contenidos"),
    @(5,  "Mantenimientos",      "SN"),
    @(6,  "Reparación",          "Manual de reparación"),
    @(7,  "Evaluación Unidad",   "Subir Portafolio"),
    @(8,  "Evaluación",          "SN"),
    @(9,  "Encuadre",            "Firmar Encuadre"),
    @(10, "Prueba Diagnóstico",  "SN"),
    @(11, "Introducción",        "SN"),
    @(12, "Contenidos Varios",   "Tipos de contenidos"),
    @(13, "Mantenimientos",      "SN"),
    @(14, "Reparación",          "Manual de reparación"),
    @(15, "Evaluación Unidad",   "Subir Portafolio"),
    @(16, "Encuadre",            "Firmar Encuadre"),
    @(17, "Prueba Diagnóstico",  "SN"),
    @(18, "Encuadre",            "Firmar Encuadre")
)

$firstRow = 2

for ($i = 0; $i -lt $actividades.Count; $i++) {
    $row = $firstRow + $i
    $fila = $actividades[$i]

    # Orden de escritura A, B, D, C, E, F, G para que las cadenas nuevas
    # se registren en ese mismo orden dentro de sharedStrings.
    $ws.Range("A$row").Value = $docente
    $ws.Range("B$row").Value = $materia
    $ws.Range("D$row").Value = $semestre
    $ws.Range("C$row").Value = $carrera
    $ws.Range("E$row").Value = $fila[0]
    $ws.Range("F$row").Value = $fila[1]
    $ws.Range("G$row").Value = $fila[2]
}

$ws.Range("I14").Select() | Out-Null
